$d = $word.ActiveDocument

function Get-ParaIndexForRange($doc, $rng) {
    $s = $rng.Start
    $e = $rng.End
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $pp = $doc.Paragraphs.Item($i)
        if ($pp.Range.Start -le $s -and $pp.Range.End -ge $e) {
            return $i
        }
    }
    return -1
}

function Find-ParaIndex($doc, $needle) {
    $rng = $doc.Content
    $ok = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        return -1
    }
    $rng.Expand(4) | Out-Null
    return Get-ParaIndexForRange $doc $rng
}

# ---------------------------------------------------------------------------
# 1) "Kaue" / "Caponero" name line: merge "Kaue" + " " runs into a single
#    "Kaue " run (no spell-check wrapper) while "Caponero" keeps its own
#    spell-check wrapped run, and the trailing " - RM 96466" run is kept.
# ---------------------------------------------------------------------------
$kaueIdx = Find-ParaIndex $d "Kaue Caponero - RM 96466"
if ($kaueIdx -gt 0) {
    $r = $d.Paragraphs.Item($kaueIdx).Range
    $xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:ind w:firstLine="566"/><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Kaue </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Caponero</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> - RM 96466</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $r.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 2) Remove the stray "Entendido, aqui esta a resposta revisada..." block
#    (an empty paragraph, the sentence paragraph, and a following empty
#    paragraph) that used to sit between "2023" and "1. Escolha do Tipo de
#    Nuvem...". Then mark "1. Escolha..." as where the page now renders its
#    break (w:lastRenderedPageBreak before its text).
# ---------------------------------------------------------------------------
$yearIdx = Find-ParaIndex $d "2023"
$escIdx = Find-ParaIndex $d "1. Escolha do Tipo de Nuvem"
if ($yearIdx -gt 0 -and $escIdx -gt $yearIdx) {
    $gapStart = $d.Paragraphs.Item($yearIdx).Range.End
    $gapEnd = $d.Paragraphs.Item($escIdx).Range.Start
    if ($gapEnd -gt $gapStart) {
        $gapRange = $d.Range($gapStart, $gapEnd)
        $gapRange.Delete()
    }
}

$escIdx2 = Find-ParaIndex $d "1. Escolha do Tipo de Nuvem para a Solução Proposta:"
if ($escIdx2 -gt 0) {
    $r = $d.Paragraphs.Item($escIdx2).Range
    $xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">1. Escolha do Tipo de Nuvem para a Solução Proposta:</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $r.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 3) IA / cloud computing paragraph: drop the ", previsão de tendências"
#    phrase and move the w:lastRenderedPageBreak marker off this run...
# ---------------------------------------------------------------------------
$iaIdx = Find-ParaIndex $d "Inteligência Artificial (IA)"
if ($iaIdx -gt 0) {
    $r = $d.Paragraphs.Item($iaIdx).Range
    $xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">- Inteligência Artificial (IA): A IA é mencionada como um pilar fundamental para aumentar a precisão e eficácia dos processos de compras. No contexto de cloud </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>computing</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>, a IA pode ser utilizada para análise de dados e suporte à tomada de decisões estratégicas, o que pode ser implementado como um serviço na nuvem.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $r.InsertXML($xml)
}

# ...and land it instead on the following "Serviços de Armazenamento..." run.
$servIdx = Find-ParaIndex $d "Serviços de Armazenamento e Processamento de Dados"
if ($servIdx -gt 0) {
    $r = $d.Paragraphs.Item($servIdx).Range
    $xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">   - Serviços de Armazenamento e Processamento de Dados: Dada a necessidade de análise de dados e comunicação eficiente com fornecedores, serviços de armazenamento e processamento de dados na nuvem podem ser fundamentais. Eles oferecem a flexibilidade e a escalabilidade necessárias para gerenciar grandes volumes de dados e realizar análises complexas.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $r.InsertXML($xml)
}
